$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Date placeholder fields ("datetimeFigureOut") across the slide
#    master, every slide layout, and the notes master: 10/29/2009 ->
#    11/5/2009.
# ------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "10/29/2009") {
                $tr.Text = "11/5/2009"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DateShape $notesMaster.Shapes

# ------------------------------------------------------------------
# 2) Title slide subtitle: "October 27, 2009" -> "November 5, 2009"
#    (split across two runs, matching how the date was retyped).
# ------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf("October 27, 2009")
if ($idx -ge 0) {
    $prefix = $tr.Characters($idx + 1, 12)   # "October 27, "
    $prefix.Text = "November 5, "
}

# ------------------------------------------------------------------
# 3) "Questions for the Client" slide: merge the two runs that make
#    up "Individual descriptions for table rows" into one run.
# ------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$content = $slide11.Shapes.Item(2)
$tr11 = $content.TextFrame.TextRange
$full11 = $tr11.Text
$target = "Individual descriptions for table rows"
$idx11 = $full11.IndexOf("Individual descriptions for table ")
if ($idx11 -ge 0) {
    $whole = $tr11.Characters($idx11 + 1, $target.Length)
    $whole.Text = $target
}
